# Update metadata in the "Informe-03-030008-A-TC-TP" sheet to reflect the
# 2016-04-06 data refresh: dimension/measure vocabulary renamed from
# "iaest-measure:" to "iaest-dimension:" for dimension columns, the four
# residencia/nacimiento "-nombre" columns consolidated under the single
# sdmx-dimension:refArea concept, datatypes updated (skos:Concept / URI-*),
# and a new row 6 added with mapping file references.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: sdmx/iaest concept identifiers ---------------------------------
$ws.Range("A3").Value = "iaest-dimension:edad-grandes-grupos"
$ws.Range("F3").Value = "sdmx-dimension:refArea"
$ws.Range("G3").Value = "sdmx-dimension:refArea"
$ws.Range("H3").Value = "sdmx-dimension:refArea"
$ws.Range("I3").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "iaest-dimension:sexo"
$ws.Range("L3").Value = "iaest-dimension:relacion-lugar-de-residencia-y-nacimiento"

# --- Row 4: measure/dim classification --------------------------------------
$ws.Range("A4").Value = "dim"
$ws.Range("F4").Value = "dim"
$ws.Range("G4").Value = "dim"
$ws.Range("H4").Value = "dim"
$ws.Range("I4").Value = "dim"
$ws.Range("K4").Value = "dim"
$ws.Range("L4").Value = "dim"

# --- Row 5: datatypes --------------------------------------------------------
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("F5").Value = "URI-Comunidad"
$ws.Range("G5").Value = "URI-Comunidad"
$ws.Range("H5").Value = "URI-comarca"
$ws.Range("I5").Value = "URI-Provincia"
$ws.Range("K5").Value = "skos:Concept"
$ws.Range("L5").Value = "skos:Concept"

# --- Row 6 (new): mapping workbook references --------------------------------
# Use Copy to bring over the same cell style ("s=1") used throughout the
# existing data rows before overwriting the value.
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = "mapping-edad-grandes-grupos.xlsx"

$ws.Range("K5").Copy($ws.Range("K6"))
$ws.Range("K6").Value = "mapping-sexo.xlsx"

$ws.Range("L5").Copy($ws.Range("L6"))
$ws.Range("L6").Value = "mapping-relacion-lugar-de-residencia-y-nacimiento.xlsx"
